# ---------------------------------------------------------------------------
# carpg/doc/stats.xlsx update
#   "debug draw stamina, training, enemies use stamina, updated texts"
#
#   - add a new worksheet "Arkusz6" with the stamina / hp debug-draw tables
#   - make it the active / selected sheet (was Arkusz5)
#   - sheet4 ("Arkusz4") selection moves from F15 to B11
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- sheet4: selection moved from F15 to B11 -------------------------------
$ws4 = $wb.Worksheets.Item("Arkusz4")
$ws4.Activate() | Out-Null
$ws4.Range("B11").Select() | Out-Null

# --- new sheet "Arkusz6", inserted after the last existing sheet ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws6.Name = "Arkusz6"

# Row 1: level header + the 5 sampled levels (0,5,10,15,20)
$ws6.Range("A1").Value = "level"
$ws6.Range("B1").Value = 0
$ws6.Range("C1").Value = 5
$ws6.Range("D1").Value = 10
$ws6.Range("E1").Value = 15
$ws6.Range("F1").Value = 20

# Row labels reusing / introducing shared strings in the order the
# original workbook grew them in
$ws6.Range("A2").Value = "str"
$ws6.Range("A3").Value = "end"
$ws6.Range("A4").Value = "dex"

# Base / bonus attribute table (rows 9-11)
$ws6.Range("A9").Value = "base str"
$ws6.Range("B9").Value = 65
$ws6.Range("A10").Value = "base end"
$ws6.Range("B10").Value = 65
$ws6.Range("A11").Value = "base dex"
$ws6.Range("B11").Value = 55

$ws6.Range("C10").Value = "bonus end"
$ws6.Range("D10").Value = 6.25
$ws6.Range("C9").Value = "bonus str"
$ws6.Range("D9").Value = 6.25
$ws6.Range("C11").Value = "bonus dex"
$ws6.Range("D11").Value = 3.75

$ws6.Range("A6").Value = "hp"
$ws6.Range("A7").Value = "stamina"
$ws6.Range("A14").Value = "hp v"

# Row 2: str progression = base str + bonus str * level/5
$ws6.Range("B2").Formula = '=$B9+$D9*B$1/5'
$ws6.Range("C2:F2").Formula = '=$B9+$D9*C$1/5'

# Row 3: end progression = base end + bonus end * level/5
$ws6.Range("B3:F3").Formula = '=$B10+$D10*B$1/5'

# Row 4: dex progression = base dex + bonus dex * level/5
$ws6.Range("B4:F4").Formula = '=$B11+$D11*B$1/5'

# Row 6: hp = 500 * (1 + (hp v - 50)/50)
$ws6.Range("B6").Formula = '=500 * (1 + (B14-50)/50)'
$ws6.Range("C6:E6").Formula = '=500 * (1 + (C14-50)/50)'
$ws6.Range("F6").Formula = '=500 * (1 + (F14-50)/50)'

# Row 7: stamina = 50 + end*2.5 + dex*2
$ws6.Range("B7").Formula = '=50+B3*2.5+B4*2'
$ws6.Range("C7:F7").Formula = '=50+C3*2.5+C4*2'

# Row 14: hp v = end*0.8 + str*0.2
$ws6.Range("B14").Formula = '=B3*0.8 + B2*0.2'
$ws6.Range("C14:F14").Formula = '=C3*0.8 + C2*0.2'

# --- Arkusz6 becomes the active / selected sheet --------------------------
$ws6.Activate() | Out-Null
$ws6.Range("I10").Select() | Out-Null
